# Update gh-pages output (江西-漫展信息.xlsx) to reflect the newly
# scraped convention-listing data: a new event ("信州·上漫·ACG动漫游戏
# 嘉年华") was inserted as row 3 on both the "展览" and "全部类型" sheets
# (pushing the later rows down by one), and a handful of "想去人数"
# (interest-count) numbers were refreshed.

$wb = $excel.ActiveWorkbook

function Insert-NewEvent($ws, $RowIndex, $SeqNo) {

    # Push everything from $RowIndex down by one row, then populate the
    # freshly-opened row with the new event's data.
    $ws.Rows.Item($RowIndex).Insert()

    $row = $ws.Rows.Item($RowIndex)

    # Column A is the 0-based sequence number; reuse the same bold /
    # centered / bordered look as every other row in column A instead of
    # whatever formatting Insert() happened to carry down.
    $aCell = $ws.Cells.Item($RowIndex, 1)
    $aCell.Value = $SeqNo
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108   # xlCenter
    $aCell.VerticalAlignment = -4160     # xlTop
    $aCell.Borders.LineStyle = 1         # xlContinuous

    # Column B holds a plain "yyyy-MM-dd" string in this workbook, not a
    # real date -- force text so Excel doesn't auto-convert it.
    $bCell = $ws.Cells.Item($RowIndex, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = "2024-11-02"
    $bCell.Style = "Normal"

    $ws.Cells.Item($RowIndex, 3).Value = "信州·上漫·ACG动漫游戏嘉年华"
    $ws.Cells.Item($RowIndex, 4).Value = "高铁经济试验区吴楚大道与凤凰东大道交叉口 饶派沉浸式街区(B1)"
    $ws.Cells.Item($RowIndex, 5).Value = "2024.11.02 10:00-11.02 17:30"
    $ws.Cells.Item($RowIndex, 6).Value = 0
    $ws.Cells.Item($RowIndex, 7).Value = 39.9
    $ws.Cells.Item($RowIndex, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93223"
    $ws.Cells.Item($RowIndex, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/5BN1Xhzp1728294588740.jpeg"
}

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) -- before: rows 2-10, after: rows 2-11
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 746

Insert-NewEvent $ws1 3 2

$ws1.Range("F4").Value = 43      # 南昌·花绒万兽秋镜派对
$ws1.Range("F6").Value = 3197    # 南昌·CM04动漫游戏博览会
$ws1.Range("F8").Value = 3930    # 南昌·云芽动漫音乐嘉年华
$ws1.Range("F9").Value = 484     # 南昌·云芽动漫音乐嘉年华·封茗囧菌内场票
$ws1.Range("F10").Value = 989    # 南昌·萌卡动漫展
$ws1.Range("F11").Value = 36     # 九江·第二届异次元动漫嘉年华

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) -- before: rows 2-11, after: rows 2-12
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 746

Insert-NewEvent $ws4 3 2

$ws4.Range("F4").Value = 43      # 南昌·花绒万兽秋镜派对
$ws4.Range("F7").Value = 3197    # 南昌·CM04动漫游戏博览会
$ws4.Range("F9").Value = 3930    # 南昌·云芽动漫音乐嘉年华
$ws4.Range("F10").Value = 484    # 南昌·云芽动漫音乐嘉年华·封茗囧菌内场票
$ws4.Range("F11").Value = 989    # 南昌·萌卡动漫展
$ws4.Range("F12").Value = 36     # 九江·第二届异次元动漫嘉年华

Write-Host "Event list refreshed."
